$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: paste formats only (xlPasteFormats = -4122) from a template cell
# ---------------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 1 - add a "NER" label in the new column P (header row)
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "NER"

# ---------------------------------------------------------------------------
# Row 9 - continuation of the bert_en_uncased / NER(RONEC) table
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = "bert_en_uncased"
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 256
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 0.8834
$ws.Range("F9").Value = 0.6878
$ws.Range("G9").Value = 0.6349
$ws.Range("H9").Value = 0.656
$ws.Range("I9").Value = 0.4577
$ws.Range("J9").Value = 11178785
$ws.Range("K9").Value = "Adam with polynomial decay (init_lr=2e-5, end_lr=1e-9)"
$ws.Range("L9").Value = 32
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = "RONEC"

Copy-Format "E2" "E9"
Copy-Format "E2" "F9"
Copy-Format "E2" "G9"
Copy-Format "E2" "H9"
Copy-Format "E2" "I9"
Copy-Format "J3" "J9"
Copy-Format "E2" "L9"
Copy-Format "E2" "M9"

# ---------------------------------------------------------------------------
# Row 19 - header row for the new "Units / Embedding_output_dim" table
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Name"
$ws.Range("B19").Value = "L"
$ws.Range("C19").Value = "Units"
$ws.Range("D19").Value = "Embedding_output_dim"
$ws.Range("E19").Value = "Accuracy"
$ws.Range("F19").Value = "Precision"
$ws.Range("G19").Value = "Recall"
$ws.Range("H19").Value = "F1"
$ws.Range("I19").Value = "Loss"
$ws.Range("J19").Value = "Parameters"
$ws.Range("K19").Value = "Optimizer"
$ws.Range("L19").Value = "Batch size"
$ws.Range("N19").Value = "Dataset"
$ws.Range("P19").Value = "Text classification"

# ---------------------------------------------------------------------------
# Rows 20-22 - rnn_bidirectional_en results
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "rnn_bidirectional_en"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = 64
$ws.Range("E20").Value = 0.9205
$ws.Range("F20").Value = 0.9151
$ws.Range("G20").Value = 0.9151
$ws.Range("H20").Value = 0.9001
$ws.Range("I20").Value = 1.5924
$ws.Range("J20").Value = 34254
$ws.Range("K20").Value = "Adam(lr=1e-4)"
$ws.Range("L20").Value = 12
$ws.Range("N20").Value = "en_text_classification + en_personas"

$ws.Range("A21").Value = "rnn_bidirectional_en"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = 0.9315
$ws.Range("F21").Value = 0.9202
$ws.Range("G21").Value = 0.9202
$ws.Range("H21").Value = 0.9126
$ws.Range("I21").Value = 1.2462
$ws.Range("J21").Value = 76366
$ws.Range("K21").Value = "Adam(lr=1e-4)"
$ws.Range("L21").Value = 12
$ws.Range("N21").Value = "en_text_classification + en_personas"

$ws.Range("A22").Value = "rnn_bidirectional_en"
$ws.Range("B22").Value = 2
$ws.Range("C22").Value = 64
$ws.Range("D22").Value = 64
$ws.Range("E22").Value = 0.6438
$ws.Range("F22").Value = 0.7206
$ws.Range("G22").Value = 0.7206
$ws.Range("H22").Value = 0.7405
$ws.Range("I22").Value = 0.9991
$ws.Range("J22").Value = 175182
$ws.Range("K22").Value = "Adam(lr=1e-4)"
$ws.Range("L22").Value = 12
$ws.Range("N22").Value = "en_text_classification + en_personas"

"E20","F20","G20","H20","I20","E21","F21","G21","H21","I21","E22","F22","G22","H22","I22" | ForEach-Object {
    Copy-Format "E2" $_
}
"J20","J21","J22" | ForEach-Object {
    Copy-Format "J3" $_
}

# ---------------------------------------------------------------------------
# Row 28 - header row for the new bert_en_uncased (English) table
# ---------------------------------------------------------------------------
$ws.Range("A28").Value = "Name"
$ws.Range("B28").Value = "L"
$ws.Range("C28").Value = "H"
$ws.Range("D28").Value = "A"
$ws.Range("E28").Value = "Accuracy"
$ws.Range("F28").Value = "Precision"
$ws.Range("G28").Value = "Recall"
$ws.Range("H28").Value = "F1"
$ws.Range("I28").Value = "Loss"
$ws.Range("J28").Value = "Parameters"
$ws.Range("K28").Value = "Optimizer"
$ws.Range("L28").Value = "Batch size"
$ws.Range("M28").Value = "Sequence length"
$ws.Range("N28").Value = "Dataset"
$ws.Range("P28").Value = "Text classification"

# ---------------------------------------------------------------------------
# Row 29 - bert_en_uncased English text classification result
# ---------------------------------------------------------------------------
$ws.Range("A29").Value = "bert_en_uncased"
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = 256
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = 0.9717
$ws.Range("F29").Value = 0.9741
$ws.Range("G29").Value = 0.9741
$ws.Range("H29").Value = 0.964
$ws.Range("I29").Value = 0.2371
$ws.Range("J29").Value = 11174159
$ws.Range("K29").Value = "Adam(lr=1e-5)"
$ws.Range("L29").Value = 12
$ws.Range("M29").Value = 128
$ws.Range("N29").Value = "English"

"E29","F29","G29","H29","I29","L29","M29" | ForEach-Object {
    Copy-Format "E2" $_
}
Copy-Format "J3" "J29"

# ---------------------------------------------------------------------------
# Column widths (the underlying engine rounds ColumnWidth to whole pixels,
# same as real Excel - inputs below are tuned to land as close as possible
# on the target stored widths from the source workbook)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.5847
$ws.Columns.Item(4).ColumnWidth = 21.2588
$ws.Columns.Item(9).ColumnWidth = 10.5847
$ws.Columns.Item(10).ColumnWidth = 12.9173
$ws.Columns.Item(14).ColumnWidth = 32.422
$ws.Columns.Item(16).ColumnWidth = 16.422

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("N6").Select() | Out-Null
